# Generate Report for Handoff
# Adds a new handed-off file ("756b0143-fbe4-476a-be2a-a9c675e87d20.md") as a
# third row to the Overview sheet and to each locale sheet (zh-cn, de-de),
# growing each table by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) -> table "Overview", range A1:G2 -> A1:G3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "756b0143-fbe4-476a-be2a-a9c675e87d20.md"

$wsOverview.Range("B3").Value = "e2e\756b0143-fbe4-476a-be2a-a9c675e87d20.md"
$wsOverview.Range("B3").Style = "HyperLink"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8925e133915d25e0ba8cc8b45e8bd99d9863b5a5/e2e/756b0143-fbe4-476a-be2a-a9c675e87d20.md", [Type]::Missing, [Type]::Missing, "e2e\756b0143-fbe4-476a-be2a-a9c675e87d20.md") | Out-Null

$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"

$wsOverview.Range("G3").Value = "2016-08-22 02:50:35"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) -> table "zh_cn", range A1:P2 -> A1:P3
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = "756b0143-fbe4-476a-be2a-a9c675e87d20.md"
$wsZhCn.Range("A3").Style = "HyperLink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8925e133915d25e0ba8cc8b45e8bd99d9863b5a5/e2e/756b0143-fbe4-476a-be2a-a9c675e87d20.md", [Type]::Missing, [Type]::Missing, "756b0143-fbe4-476a-be2a-a9c675e87d20.md") | Out-Null

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "756b0143-fbe4-476a-be2a-a9c675e87d20.da91853ba9c461955e9e50afee3ffd9fbb7b0b46.zh-cn.xlf"

$wsZhCn.Range("H3").Value = "2016-08-22 02:50:31"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""

$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3) -> table "de_de", range A1:P2 -> A1:P3
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = "756b0143-fbe4-476a-be2a-a9c675e87d20.md"
$wsDeDe.Range("A3").Style = "HyperLink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8925e133915d25e0ba8cc8b45e8bd99d9863b5a5/e2e/756b0143-fbe4-476a-be2a-a9c675e87d20.md", [Type]::Missing, [Type]::Missing, "756b0143-fbe4-476a-be2a-a9c675e87d20.md") | Out-Null

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "756b0143-fbe4-476a-be2a-a9c675e87d20.da91853ba9c461955e9e50afee3ffd9fbb7b0b46.de-de.xlf"

$wsDeDe.Range("H3").Value = "2016-08-22 02:50:35"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""

$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""
